# "search delete view testing"
# Delete the value in G1 (student_college_id header) and update the
# active selection from A2:G4 down to just A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "student_college_id" header value from G1, leaving the cell
# (and its formatting/style) in place but empty.
$ws.Range("G1").ClearContents()

# Update the worksheet selection to just A2 (was A2:G4, active cell G4).
$ws.Range("A2").Select()
